$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-12 (columns D, L, M, N, O, P, Q, R, S, T).
# Rows 2-12 all share identical A,B,C,E,F,G,H,I,J,K values, so this is a
# re-shuffle of the variable columns between existing rows.
$rows = @(
    @{ Row = 2;  D = 45086; L = "Primera"; M = 30;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";      S = 1000; T = 18 },
    @{ Row = 3;  D = 44698; L = "Primera"; M = 120; N = 16000; O = 17000; P = 16500; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins";   S = 917;  T = 18 },
    @{ Row = 4;  D = 45099; L = "Primera"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";      S = 1000; T = 18 },
    @{ Row = 5;  D = 44344; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó";   S = 750;  T = 18 },
    @{ Row = 6;  D = 45085; L = "Primera"; M = 60;  N = 18000; O = 19000; P = 18500; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";      S = 1028; T = 18 },
    @{ Row = 7;  D = 44334; L = "Primera"; M = 120; N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada"; R = "Región de O'Higgins"; S = 1042; T = 12 },
    @{ Row = 8;  D = 44316; L = "Primera"; M = 60;  N = 17500; O = 18000; P = 17750; Q = "`$/caja 16 kilos granel"; R = "Región de O'Higgins";   S = 1109; T = 16 },
    @{ Row = 9;  D = 44316; L = "Segunda"; M = 40;  N = 16000; O = 16000; P = 16000; Q = "`$/caja 16 kilos granel"; R = "Región de O'Higgins";   S = 1000; T = 16 },
    @{ Row = 10; D = 45092; L = "Primera"; M = 35;  N = 18000; O = 19000; P = 18571; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";      S = 1032; T = 18 },
    @{ Row = 11; D = 44330; L = "Primera"; M = 60;  N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó";   S = 861;  T = 18 },
    @{ Row = 12; D = 45096; L = "Primera"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";      S = 1000; T = 18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N: Precio mínimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Precio máximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Unidad de comercialización
    $ws.Cells.Item($row, 18).Value = $r.R   # R: Origen
    $ws.Cells.Item($row, 19).Value = $r.S   # S: Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $r.T   # T: Kg / unidad
}
